$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.614.18"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.881.67"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5119"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3949"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08427"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.294"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.880.85"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.296"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001110"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06736"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.976"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "28.656.62"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.255"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "2.101.35"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.388"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1055"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.832"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.625"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02470"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06547"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2196"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.955"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.263"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.105"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6479"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.007"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6089"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.708"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.049"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.195"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.78%  "
